$wb = $excel.ActiveWorkbook

# Sheet 2: "Alias Datatype Usage Proper" - delete row 20 entirely (shifts rows 21+ up by one)
$ws2 = $wb.Worksheets.Item("Alias Datatype Usage Proper")
[void]$ws2.Activate()
[void]$ws2.Rows.Item(20).Delete()

# Update selection on sheet2 to F31 (post-shift position), matching the new
# active selection recorded after the edit.
[void]$ws2.Range("F31").Select()
